# Update the Haba (fava bean) weekly price data.
# Each data row (3-13) is refreshed with the values coming from a different
# row of the previous snapshot (rows rotate as new weekly data arrives),
# as captured by the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: target row -> source row (values to copy into target row)
$rowMap = @{
    3  = 13
    4  = 12
    5  = 9
    6  = 10
    7  = 3
    8  = 4
    9  = 5
    10 = 11
    11 = 7
    12 = 6
    13 = 8
}

# Capture original values (D, J, K, L, M, O, P) for rows 3-13 before any
# cell is overwritten, so the row rotation is computed from the original
# snapshot rather than partially-updated data. Value2 is used so dates
# come back as raw serial numbers instead of formatted text.
$orig = @{}
foreach ($r in 3..13) {
    $orig[$r] = @{
        D = $ws.Range("D$r").Value2
        J = $ws.Range("J$r").Value2
        K = $ws.Range("K$r").Value2
        L = $ws.Range("L$r").Value2
        M = $ws.Range("M$r").Value2
        O = $ws.Range("O$r").Value2
        P = $ws.Range("P$r").Value2
    }
}

foreach ($r in 3..13) {
    $src = $rowMap[$r]
    $vals = $orig[$src]
    $ws.Range("D$r").Value2 = $vals.D
    $ws.Range("J$r").Value2 = $vals.J
    $ws.Range("K$r").Value2 = $vals.K
    $ws.Range("L$r").Value2 = $vals.L
    $ws.Range("M$r").Value2 = $vals.M
    $ws.Range("O$r").Value2 = $vals.O
    $ws.Range("P$r").Value2 = $vals.P
}
